$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Price (column D) tweaks -- stored as text in the sheet, so force text
# storage via NumberFormat "@" to avoid numeric auto-coercion.
Set-TextValue "D2"  "266.07"
Set-TextValue "D3"  "21.30"
Set-TextValue "D4"  "6.114"
Set-TextValue "D5"  "0.06100"
Set-TextValue "D7"  "6.490"
Set-TextValue "D8"  "1.363"
Set-TextValue "D9"  "0.8213"
Set-TextValue "D11" "0.1587"
Set-TextValue "D12" "0.08102"
Set-TextValue "D13" "0.03408"
Set-TextValue "D14" "0.03204"
Set-TextValue "D16" "3.745"
Set-TextValue "D17" "0.001630"
Set-TextValue "D18" "0.04650"
Set-TextValue "D19" "0.006387"
Set-TextValue "D20" "0.006145"
Set-TextValue "D21" "0.001069"
Set-TextValue "D22" "0.0001501"
Set-TextValue "D23" "3.725"
Set-TextValue "D24" "2.268"
Set-TextValue "D40" "0.04592"
Set-TextValue "D41" "0.006991"

# Rows 42/43 swap coin identity (BKEXToken <-> CEJI) along with new prices
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003902"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1115"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue "D44" "0.01169"
Set-TextValue "D45" "0.00005961"
Set-TextValue "D48" "0.8024"
Set-TextValue "D50" "0.00001901"
Set-TextValue "D51" "0.01241"
